$wb = $excel.ActiveWorkbook

# --- Add the new "Facebook" worksheet after the existing "Merge Intervals" sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Facebook"

# --- Column E gets the custom date format applied at the column level first ---
$ws.Columns("E").NumberFormat = "[$-F800]dddd, mmmm dd, yyyy"

# --- Header row ---
$ws.Range("A1").Value = "Question Url"
$ws.Range("B1").Value = "Solution"
$ws.Range("C1").Value = "Solved?"
$ws.Range("D1").Value = "How Long"
$ws.Range("E1").Value = "First Solution"
$ws.Range("F1").Value = "Date of Review 1"
$ws.Range("G1").Value = "Date of Review 2"
$ws.Range("H1").Value = "Date of Review 3"

$ws.Range("A1:D1,F1:H1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108

# --- Data row ---
$url = "https://leetcode.com/problems/valid-palindrome-ii/"
$ws.Hyperlinks.Add($ws.Range("A2"), $url, "", "", $url) | Out-Null
$ws.Range("A2").Value = "Valid Palindrome II - LeetCode"
$ws.Range("A2").Style = "Hyperlink"

$ws.Range("B2").Value = "2 Pointers, Loop, check if left is not equal to right, if a character is deleted, return false, else, try deleting a character by recursing through the function"
$ws.Range("B2").WrapText = $true

$ws.Range("C2").Value = "No"
$ws.Range("C2").HorizontalAlignment = -4108

$ws.Range("D2").Value = "Looked at discussions"

$ws.Range("E2").Value = 44391

$ws.Rows(2).RowHeight = 60

# --- Column widths (character units, matching the authored workbook) ---
$ws.Columns("A").ColumnWidth = 30.666666666666664
$ws.Columns("B").ColumnWidth = 47.16666666666667
$ws.Columns("C").ColumnWidth = 7.166666666666666
$ws.Columns("D").ColumnWidth = 19.5
$ws.Columns("E").ColumnWidth = 19.5
$ws.Columns("F:H").ColumnWidth = 15.166666666666668

# --- Page setup ---
$ws.PageSetup.Orientation = 1

$wb.Save()
